$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing rows for the MuSCs sending-cluster combinations
# that are no longer present after the TPM data refresh
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf2"
$ws.Range("C2").Value = "Fgfr4"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1817723333333333
$ws.Range("N2").Value = 0.5453170000000001
$ws.Range("O2").Value = 0.008966262009224884
$ws.Range("P2").Value = 0.008966262009224884
$ws.Range("Q2").Value = 0.02975037484277778
$ws.Range("R2").Value = 0.2677533735850001
$ws.Range("S2").Value = 0.00007746395837868003
$ws.Range("T2").Value = 0.00007746395837868005

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf2"
$ws.Range("C3").Value = "Fgfr4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.09115
$ws.Range("N3").Value = 60.27345
$ws.Range("O3").Value = 0.9910337379907751
$ws.Range("P3").Value = 0.9910337379907752
$ws.Range("Q3").Value = 3.28828503525
$ws.Range("R3").Value = 29.59456531725
$ws.Range("S3").Value = 0.008562029098926774
$ws.Range("T3").Value = 0.008562029098926775

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Fgf2"
$ws.Range("C4").Value = "Fgfr4"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 15.322826
$ws.Range("H4").Value = 45.968478
$ws.Range("I4").Value = 0.808839719627903
$ws.Range("J4").Value = 0.8088397196279031
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1817723333333333
$ws.Range("N4").Value = 0.5453170000000001
$ws.Range("O4").Value = 0.008966262009224884
$ws.Range("P4").Value = 0.008966262009224884
$ws.Range("Q4").Value = 2.785265835280667
$ws.Range("R4").Value = 25.06739251752601
$ws.Range("S4").Value = 0.007252268849651773
$ws.Range("T4").Value = 0.007252268849651774

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf2"
$ws.Range("C5").Value = "Fgfr4"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 15.322826
$ws.Range("H5").Value = 45.968478
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.09115
$ws.Range("N5").Value = 60.27345
$ws.Range("O5").Value = 0.9910337379907751
$ws.Range("P5").Value = 0.9910337379907752
$ws.Range("Q5").Value = 307.8531955899
$ws.Range("R5").Value = 2770.6787603091
$ws.Range("S5").Value = 0.8015874507782512
$ws.Range("T5").Value = 0.8015874507782514

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Fgf2"
$ws.Range("C6").Value = "Fgfr4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.457711333333334
$ws.Range("H6").Value = 10.373134
$ws.Range("I6").Value = 0.1825207873147914
$ws.Range("J6").Value = 0.1825207873147914
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1817723333333333
$ws.Range("N6").Value = 0.5453170000000001
$ws.Range("O6").Value = 0.008966262009224884
$ws.Range("P6").Value = 0.008966262009224884
$ws.Range("Q6").Value = 0.6285162570531112
$ws.Range("R6").Value = 5.656646313478
$ws.Range("S6").Value = 0.001636529201194429
$ws.Range("T6").Value = 0.001636529201194429

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Fgf2"
$ws.Range("C7").Value = "Fgfr4"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.457711333333334
$ws.Range("H7").Value = 10.373134
$ws.Range("I7").Value = 0.1825207873147914
$ws.Range("J7").Value = 0.1825207873147914
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 20.09115
$ws.Range("N7").Value = 60.27345
$ws.Range("O7").Value = 0.9910337379907751
$ws.Range("P7").Value = 0.9910337379907752
$ws.Range("Q7").Value = 69.4693970547
$ws.Range("R7").Value = 625.2245734923
$ws.Range("S7").Value = 0.180884258113597
$ws.Range("T7").Value = 0.180884258113597
